$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
for ($i=1; $i -le $ws.Shapes.Count; $i++) {
    $sh = $ws.Shapes.Item($i)
    Write-Host $i $sh.Name $sh.Top $sh.Left
}
